$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.510.86'
$ws.Range("E2").Value = '  +3.60%  '
$ws.Range("D3").Value = '2.073.38'
$ws.Range("E3").Value = '  +3.16%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.651'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.84'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.99%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.396'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0817'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.15%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.925'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +20.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '2.377.36'
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.66'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.93%  '
$ws.Range("D18").Value = '2.069.41'
$ws.Range("E18").Value = '  +3.22%  '
$ws.Range("D19").Value = '37.466.06'
$ws.Range("E19").Value = '  +3.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.85%  '
$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").Value = '  +5.35%  '
$ws.Range("E22").Value = '  +4.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.34%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.31%  '
$ws.Range("E26").Value = '  +3.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.05'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.128'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +28.66%  '
$ws.Range("E31").Value = '  +2.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0630'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("E36").Value = '  +2.45%  '
$ws.Range("E37").Value = '  +11.08%  '
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +30.47%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.29'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.42%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.102'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.19%  '
$ws.Range("E43").Value = '  +4.04%  '
$ws.Range("E44").Value = '  +4.53%  '
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.34%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("D49").Value = '1.399.13'
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.17%  '
